$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.867.15'
$ws.Range('E2').Value = '  -4.10%  '
$ws.Range('D3').Value = '1.956.27'
$ws.Range('E3').Value = '  -4.13%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '241.81'
$ws.Range('E5').Value = '  -4.21%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  -4.09%  '
$ws.Range('D7').Value = '61.58'
$ws.Range('E7').Value = '  -5.44%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.364'
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('D10').Value = '56.12'
$ws.Range('E10').Value = '  -4.97%  '
$ws.Range('D11').Value = '0.0791'
$ws.Range('E11').Value = '  +4.86%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '0.854'
$ws.Range('E13').Value = '  -5.72%  '
$ws.Range('D14').Value = "'13.90"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.98%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '21.58'
$ws.Range('E15').Value = '  +3.62%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.238.73'
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('D17').Value = '5.39'
$ws.Range('E17').Value = '  -3.85%  '
$ws.Range('D18').Value = '1.962.26'
$ws.Range('E18').Value = '  -4.52%  '
$ws.Range('D19').Value = '35.704.07'
$ws.Range('E19').Value = '  -4.40%  '
$ws.Range('D20').Value = '70.51'
$ws.Range('E20').Value = '  -3.56%  '
$ws.Range('D21').Value = '0.0₃0848'
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('D22').Value = '238.54'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('D23').Value = '5.17'
$ws.Range('E23').Value = '  -3.48%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').Value = '  -8.77%  '
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  -2.69%  '
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = '158.07'
$ws.Range('E28').Value = '  -4.71%  '
$ws.Range('D29').Value = '19.66'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = "'0.130"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +14.56%  '
$ws.Range('D31').Value = '0.119'
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('D32').Value = '4.85'
$ws.Range('E32').Value = '  -7.03%  '
$ws.Range('D33').Value = '1.13'
$ws.Range('E33').Value = '  -7.46%  '
$ws.Range('D34').Value = '0.0613'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '4.36'
$ws.Range('E35').Value = '  -7.77%  '
$ws.Range('D36').Value = '6.24'
$ws.Range('E36').Value = '  +4.59%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.29'
$ws.Range('E37').Value = '  -6.66%  '
$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').Value = '3.11'
$ws.Range('E40').Value = '  +13.68%  '
$ws.Range('D41').Value = '0.0982'
$ws.Range('E41').Value = '  -5.21%  '
$ws.Range('D42').Value = '1.21'
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0211'
$ws.Range('E43').Value = '  -3.69%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -4.70%  '
$ws.Range('D45').Value = '1.08'
$ws.Range('E45').Value = '  -5.07%  '
$ws.Range('D46').Value = '92.07'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('D47').Value = '16.04'
$ws.Range('E47').Value = '  -6.31%  '
$ws.Range('D48').Value = '7.48'
$ws.Range('E48').Value = '  -8.15%  '
$ws.Range('D49').Value = '1.336.40'
$ws.Range('E49').Value = '  -5.59%  '
$ws.Range('D50').Value = '2.74'
$ws.Range('E50').Value = '  -6.48%  '
$ws.Range('D51').Value = '2.132.96'
$ws.Range('E51').Value = '  -4.26%  '
